$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column: crypto price text (some look numeric and need a text-forcing
# apostrophe prefix so Excel keeps them as literal text instead of parsing
# them into a float, matching the original inlineStr cell data).
# E column: percentage-change text, never numeric-looking (has "%" and
# surrounding spaces), so it is always safe to set directly.

$ws.Range('D2').Value = '34.401.23'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').Value = '1.789.68'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = "'226.32"
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('D6').Value = "'0.554"
$ws.Range('E6').Value = '  +1.59%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = "'32.83"
$ws.Range('E8').Value = '  +2.85%  '
$ws.Range('E9').Value = '  +0.96%  '
$ws.Range('E10').Value = '  +0.62%  '
$ws.Range('D11').Value = "'0.0946"
$ws.Range('E11').Value = '  -0.29%  '
$ws.Range('D12').Value = '2.047.56'
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('D14').Value = '1.794.89'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('D15').Value = "'0.635"
$ws.Range('E15').Value = '  +1.87%  '
$ws.Range('D16').Value = '34.362.95'
$ws.Range('E16').Value = '  +0.84%  '
$ws.Range('E17').Value = '  +2.79%  '
$ws.Range('D18').Value = "'68.55"
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('D19').Value = "'245.29"
$ws.Range('E19').Value = '  -0.12%  '
$ws.Range('D20').Value = '0.0₃0796'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('E21').Value = '  +3.39%  '
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('E23').Value = '  +1.47%  '
$ws.Range('D24').Value = "'167.58"
$ws.Range('E24').Value = '  +3.50%  '
$ws.Range('E25').Value = '  +1.30%  '
$ws.Range('D26').Value = "'7.33"
$ws.Range('E26').Value = '  +2.93%  '
$ws.Range('D27').Value = "'16.58"
$ws.Range('E27').Value = '  +1.95%  '
$ws.Range('E28').Value = '  +1.56%  '
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('D30').Value = "'4.00"
$ws.Range('E30').Value = '  +7.62%  '
$ws.Range('D31').Value = "'0.0527"
$ws.Range('E31').Value = '  +1.97%  '
$ws.Range('D32').Value = "'3.81"
$ws.Range('E32').Value = '  +2.52%  '
$ws.Range('E33').Value = '  +0.21%  '
$ws.Range('E34').Value = '  +1.61%  '
$ws.Range('D35').Value = "'2.59"
$ws.Range('E35').Value = '  +5.56%  '
$ws.Range('D36').Value = '1.410.75'
$ws.Range('E36').Value = '  -2.60%  '
$ws.Range('D37').Value = "'0.684"
$ws.Range('E37').Value = '  +5.13%  '
$ws.Range('E38').Value = '  +3.13%  '
$ws.Range('E39').Value = '  +0.35%  '
$ws.Range('D40').Value = "'84.23"
$ws.Range('E40').Value = '  +4.55%  '
$ws.Range('E42').Value = '  +0.66%  '
$ws.Range('E43').Value = '  +2.92%  '
$ws.Range('D44').Value = "'13.68"
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('E45').Value = '  +1.57%  '
$ws.Range('E46').Value = '  +2.88%  '
$ws.Range('D47').Value = "'6.07"
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').Value = '1.948.54'
$ws.Range('D49').Value = "'105.48"
$ws.Range('E49').Value = '  +0.81%  '
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('E51').Value = '  -2.62%  '
